$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header O1 from "加碼值" to "合約加碼值" (shared-string table reorders
# as a side effect: the old "加碼值" entry is dropped and the new text is
# appended at the end, shifting "檢核訊息"/"擬調利率" up by one index each;
# this matches the O1/P1/V1 <v> index changes in the diff).
$ws.Range("O1").Value = "合約加碼值"

# Column O gets a bit wider to fit the longer header text.
$ws.Columns("O").ColumnWidth = 10.88671875

# Selection moved from P5 to O6.
$ws.Range("O6").Select()
